$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column H slightly (12.71 -> 13.71)
$ws.Columns.Item(8).ColumnWidth = 13.7109375

# Update the numeric results for rows 2 and 3 (random_forest, lsboost) with new values
$ws.Range("B2").Value = 2.0600074572046672
$ws.Range("C2").Value = 0.21193492358072705
$ws.Range("D2").Value = 1.7581105820105822
$ws.Range("E2").Value = 0.91980898295015612
$ws.Range("F2").Value = 0.95906672497285406
$ws.Range("G2").Value = 0.92727351371866185
$ws.Range("H2").Value = 0.080191017049843882
$ws.Range("I2").Value = 0.64284033472311763

$ws.Range("B3").Value = 1.0776068467591551
$ws.Range("C3").Value = 0.11086490192995421
$ws.Range("D3").Value = 1.0735354418487462
$ws.Range("E3").Value = 0.25169856861934481
$ws.Range("F3").Value = 0.50169569324376784
$ws.Range("G3").Value = 0.56621067608056252
$ws.Range("H3").Value = 0.74830143138065519
$ws.Range("I3").Value = 0.94119535275768729
